$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The empty paragraph right after "{{END-FOR p2}}" (and right before the
#    "arzt" attachments heading) becomes a manual page break, matching the
#    page break already used at the end of part 1 / before part 3.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute('{{END-FOR p2}}^p', $false, $false, $false, $false, $false, `
    $true, 1, $false, '{{END-FOR p2}}^p^m', 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Replace the "arzt.attachmentsHeading" / attachments / signature-blocks
#    loops with the new "arzt.liability*" heading, loop and sign-off lines.
#    The block being replaced is 7 paragraphs; the new block is 9
#    paragraphs, so 2 extra paragraphs are appended at the end of the block.
# ---------------------------------------------------------------------------

# Paragraph indices (1-based, Word COM style) of the untouched document,
# counted from the already-applied page break above.
$pHeading   = $d.Paragraphs(44)   # was "{{arzt.attachmentsHeading}}"
$pHeading.Range.Text = '{{arzt.liabilityHeading}}'

$pBlank1    = $d.Paragraphs(45)   # was "{{FOR aItem IN arzt.attachments}}"
$pBlank1.Range.Text = ''

$pForOpen   = $d.Paragraphs(46)   # was "{{INS '- ' + $aItem}}"
$pForOpen.Range.Text = '{{FOR liability IN arzt.liabilityParagraphs}}'

$pForIns    = $d.Paragraphs(47)   # was "{{END-FOR aItem}}"
$pForIns.Range.Text = '{{INS $liability}}'

$pForClose  = $d.Paragraphs(48)   # was "{{FOR aSig IN arzt.signatureBlocks}}"
$pForClose.Range.Text = '{{END-FOR liability}}'

$pBlank2    = $d.Paragraphs(49)   # was "{{INS $aSig.label + ': ' + $aSig.name}}"
$pBlank2.Range.Text = ''

$pDate      = $d.Paragraphs(50)   # was "{{END-FOR aSig}}"
$pDate.Range.Text = 'Datum: {{arzt.liabilityDateLine}}'

# Two brand-new paragraphs for the patient-name and signature lines.
$pDate.Range.InsertParagraphAfter()
$pName = $d.Paragraphs(51)
$pName.Range.Text = 'Name Patient/in: {{arzt.liabilitySignerName}}'

$pName.Range.InsertParagraphAfter()
$pSig = $d.Paragraphs(52)
$pSig.Range.Text = 'Unterschrift: ____________________'
